$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '29.162.32'
Set-TextCell $ws.Range("E2") '  -3.17%  '
Set-TextCell $ws.Range("D3") '1.849.26'
Set-TextCell $ws.Range("E3") '  -2.17%  '
Set-TextCell $ws.Range("D4") '0.9998'
Set-TextCell $ws.Range("E4") '  -0.32%  '
Set-TextCell $ws.Range("E5") '  -4.88%  '
Set-TextCell $ws.Range("D6") '238.52'
Set-TextCell $ws.Range("E6") '  -1.73%  '
Set-TextCell $ws.Range("E7") '  -0.31%  '
Set-TextCell $ws.Range("D8") '0.3057'
Set-TextCell $ws.Range("E8") '  -3.69%  '
Set-TextCell $ws.Range("D9") '0.07485'
Set-TextCell $ws.Range("E9") '  +4.66%  '
Set-TextCell $ws.Range("D10") '23.42'
Set-TextCell $ws.Range("E10") '  -5.74%  '
Set-TextCell $ws.Range("D11") '0.08121'
Set-TextCell $ws.Range("E11") '  -2.65%  '
Set-TextCell $ws.Range("D12") '1.875.11'
Set-TextCell $ws.Range("E12") '  -1.62%  '
Set-TextCell $ws.Range("D13") '0.7259'
Set-TextCell $ws.Range("E13") '  -4.24%  '
Set-TextCell $ws.Range("D14") '5.223'
Set-TextCell $ws.Range("E14") '  -3.63%  '
Set-TextCell $ws.Range("D15") '88.77'
Set-TextCell $ws.Range("E15") '  -4.61%  '
Set-TextCell $ws.Range("D16") '29.313.52'
Set-TextCell $ws.Range("E16") '  -2.80%  '
Set-TextCell $ws.Range("D17") '5.766'
Set-TextCell $ws.Range("E17") '  -6.42%  '
Set-TextCell $ws.Range("D18") '238.46'
Set-TextCell $ws.Range("E18") '  -5.05%  '
Set-TextCell $ws.Range("E19") '  -4.00%  '
Set-TextCell $ws.Range("D20") '0.000007616'
Set-TextCell $ws.Range("E20") '  -3.17%  '
Set-TextCell $ws.Range("E21") '  -0.26%  '
Set-TextCell $ws.Range("D22") '2.126.45'
Set-TextCell $ws.Range("E22") '  -3.74%  '
Set-TextCell $ws.Range("D23") '0.9999'
Set-TextCell $ws.Range("E23") '  -0.33%  '
Set-TextCell $ws.Range("D24") '7.586'
Set-TextCell $ws.Range("E24") '  -4.55%  '
Set-TextCell $ws.Range("D25") '8.998'
Set-TextCell $ws.Range("E25") '  -3.29%  '
Set-TextCell $ws.Range("D26") '161.15'
Set-TextCell $ws.Range("E26") '  -2.15%  '
Set-TextCell $ws.Range("D27") '0.1453'
Set-TextCell $ws.Range("E27") '  -7.50%  '
Set-TextCell $ws.Range("D28") '18.07'
Set-TextCell $ws.Range("E28") '  -3.49%  '
Set-TextCell $ws.Range("D29") '1.984'
Set-TextCell $ws.Range("E29") '  -3.51%  '
Set-TextCell $ws.Range("E30") '  -5.56%  '
Set-TextCell $ws.Range("D31") '4.542'
Set-TextCell $ws.Range("E31") '  -0.78%  '
Set-TextCell $ws.Range("D32") '1.493'
Set-TextCell $ws.Range("E32") '  -3.04%  '
Set-TextCell $ws.Range("D33") '3.976'
Set-TextCell $ws.Range("E33") '  -5.28%  '
Set-TextCell $ws.Range("D34") '0.05176'
Set-TextCell $ws.Range("E34") '  -3.29%  '
Set-TextCell $ws.Range("D35") '1.187'
Set-TextCell $ws.Range("E35") '  -5.44%  '
Set-TextCell $ws.Range("D36") '1.041'
Set-TextCell $ws.Range("E36") '  +3.99%  '
Set-TextCell $ws.Range("D37") '0.7015'
Set-TextCell $ws.Range("E37") '  -9.07%  '
Set-TextCell $ws.Range("D38") '2.655'
Set-TextCell $ws.Range("E38") '  -2.69%  '
Set-TextCell $ws.Range("D39") '0.01864'
Set-TextCell $ws.Range("E39") '  -4.84%  '
Set-TextCell $ws.Range("E40") '  -3.03%  '
Set-TextCell $ws.Range("D41") '0.9326'
Set-TextCell $ws.Range("E41") '  +6.78%  '
Set-TextCell $ws.Range("D42") '6.020'
Set-TextCell $ws.Range("E42") '  -0.91%  '
Set-TextCell $ws.Range("D43") '1.076.14'
Set-TextCell $ws.Range("E43") '  -2.30%  '
Set-TextCell $ws.Range("D44") '0.4286'
Set-TextCell $ws.Range("E44") '  -6.13%  '
Set-TextCell $ws.Range("D45") '70.22'
Set-TextCell $ws.Range("E45") '  -3.14%  '
Set-TextCell $ws.Range("D46") '0.9998'
Set-TextCell $ws.Range("E46") '  -0.31%  '
Set-TextCell $ws.Range("D47") '102.63'
Set-TextCell $ws.Range("E47") '  -1.71%  '
Set-TextCell $ws.Range("D48") '2.009.41'
Set-TextCell $ws.Range("E48") '  -3.99%  '
Set-TextCell $ws.Range("D49") '1.742'
Set-TextCell $ws.Range("E49") '  -6.47%  '
Set-TextCell $ws.Range("D50") '9.156'
Set-TextCell $ws.Range("E50") '  -5.13%  '
Set-TextCell $ws.Range("D51") '7.032'
Set-TextCell $ws.Range("E51") '  -7.28%  '
